$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shorten the three city names that previously included an alternate/older
# spelling in parentheses, now that the full "City, Country" name already
# lives in column A.
$ws.Range("D330").Value = "Jeddah"
$ws.Range("D363").Value = "Krakow"
$ws.Range("D449").Value = "Lucknow"

# The sheet's default column width (covering columns A through AMK, i.e.
# 1-1025) was nudged slightly narrower.
$cols = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item(1, 1025))
$cols.EntireColumn.ColumnWidth = 13.333333333333334

# Update the window scroll position / active selection left after the edit.
$win = $excel.ActiveWindow
$win.ScrollRow = 340
$win.ScrollColumn = 1
$ws.Range("A363").Select()
